# Module 11 assignment edit script
$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $findText"
    }
    return $ok
}

# 1. Clonal succession paragraph: "waiting and once triggered" -> split with "to be activated"
Replace-Text `
    "stem cells are in a dormant state waiting and once triggered any of these stem cells could differentiate and proliferate into a large population of mature cells. These stem cells are available for the lifetime of the organism. The mature clone eventually " `
    "stem cells are in a dormant state waiting to be activated and once triggered any of these stem cells could differentiate and proliferate into a large population of mature cells. These stem cells are available for the lifetime of the organism and have a limited lifespan. The mature clone eventually "

# 2. Clonal succession paragraph: "out and a new stem-cell take over" -> "out and a new stem-cell clone take over"
Replace-Text `
    " out and a new stem-cell take over for cell production." `
    " out and a new stem-cell clone take over for cell production."

# 3. Deterministic model paragraph
Replace-Text `
    ": this model assumes that the stem cells can self-renew and differentiate into a mature cell or 1 stem-cell daughter with a given probability p (mature cell: p and stem-cell:1-p)." `
    ": this model assumes that the stem cells can self-renew and differentiate into a mature cell and a stem-cell daughter. The probability of self-renewal may not be exactly 50% depending on tissue environment and may be subject to telomere length."

# 4. Stochastic model paragraph
Replace-Text `
    ": the behavior of the outcome of differentiation is random in nature; i.e., a stem cell can generate 0,1 or 2 stem cells as daughter cells; and ca be regulated like the deterministic model by tissue environment and be influenced by telomer length." `
    ": the behavior of the outcome of differentiation is random in nature; i.e., a stem cell can generate 0, 1 or 2 stem cells as daughter cells; and can be regulated like the deterministic model by factors external to the dividing cell."

# 5. Tissue engineer paragraph - big rewrite
Replace-Text `
    "A tissue engineer, would like to rely on the deterministic model, which, compared to the clonal succession or stochastic models, under deterministic conditions, can either differentiates into a mature specialized cell; for example, for tissue repair, or a stem cell to replace itself." `
    "A tissue engineer, would like to rely on the deterministic model, which is based. on the notion that stem cells exhibit a deterministic behavior given their response to differentiation stimuli."

# 6. "After a year, Geron decided..." paragraph indentation (ind left=360 twips = 18pt)
$para = $d.Paragraphs | Where-Object { $_.Range.Text -like "After a year, Geron decided to stop*" }
if ($para) {
    $para.Range.ParagraphFormat.LeftIndent = 18
    Write-Host "Indent set on paragraph."
} else {
    Write-Host "PARAGRAPH NOT FOUND for indent."
}

# 7. "a trial to treat 8 patients" -> "injected with"
Replace-Text `
    "a trial to treat 8 patients with spinal cord injury, by injecting them with" `
    "a trial to treat 8 patients with spinal cord injury, injected with"

# 8. "$170 million with $25 million ... as a" -> "from which ... was a"
Replace-Text `
    "it has already spent `$170 million with `$25 million as a loan from the" `
    "it has already spent `$170 million from which `$25 million was a loan from the"

# 9. Big "Also" paragraph rewrite
Replace-Text `
    "Also, around the same period (slide 23 11CD), until today, there has been 4 x times less NIH funding for human embryonic research compared to non-embryonic and iPSC research. Maybe around the same time, NIH guidelines were likely being communicated. In addition, Geron funded Dr Thomson research in 1998. In 2011, 3 years after, Geron did not have yet any FDA approved stem cell therapy. Geron executive committee, probably then realized that the investment needed to continue the trial but also their stem cell research; was too steep, and could have jeopardize; maybe other more promising research. With this context; it seems expected that Geron; as a public company under the pressure of investors, took the only decision they could have made and step to pursue stem cell research altogether." `
    "Also, from 2011 until today (slide 10-11D), there has been much more NIH funding in non-embryonic and iPSC research compared to human embryonic research. In addition, Geron funded Dr Thomson research in 1998. In 2011, 3 years after, Geron did not yet have any FDA approved stem cell therapy. Geron executive committee, probably then, realized that the investment needed to continue the trial but also their stem cell research; was too steep, and could jeopardize; maybe; other more promising research. With this context; it seems expected that Geron; as a public company under the pressure of investors, took the only decision they could have financially made and decided to stop pursuing stem cell research altogether."
